$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-10-20T17:07:43"
$ws.Range("V4").Value = 43.18
$ws.Range("W4").Value = 76.61
$ws.Range("X4").Value = 40.53
$ws.Range("Y4").Value = 33.47
$ws.Range("Z4").Value = 34.66
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = -48.26
$ws.Range("U6").Value = -6.59
$ws.Range("V6").Value = -6.4
$ws.Range("W6").Value = -4.83
$ws.Range("X6").Value = -2.03
$ws.Range("Y6").Value = -1.51
$ws.Range("Z6").Value = -1.01
$ws.Range("U9").Value = 94.84
$ws.Range("V9").Value = 44.66
$ws.Range("W9").Value = 78.76000000000001
$ws.Range("X9").Value = 42.22
$ws.Range("Y9").Value = 34.8
$ws.Range("Z9").Value = 36.39
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = -48.26
$ws.Range("U11").Value = -5.88
$ws.Range("V11").Value = -4.92
$ws.Range("W11").Value = -2.68
$ws.Range("X11").Value = -0.34
$ws.Range("Y11").Value = -0.17
$ws.Range("Z11").Value = 0.73
$ws.Range("U14").Value = 94.93000000000001
$ws.Range("V14").Value = 93.01000000000001
$ws.Range("W14").Value = 78.76000000000001
$ws.Range("X14").Value = 42.22
$ws.Range("Y14").Value = 34.8
$ws.Range("Z14").Value = 36.39
$ws.Range("U16").Value = -5.79
$ws.Range("V16").Value = -4.84
$ws.Range("W16").Value = -2.68
$ws.Range("X16").Value = -0.34
$ws.Range("Y16").Value = -0.17
$ws.Range("Z16").Value = 0.73
$ws.Range("U19").Value = 93.52
$ws.Range("V19").Value = 90.93000000000001
$ws.Range("W19").Value = 76.26000000000001
$ws.Range("X19").Value = 40.34
$ws.Range("Y19").Value = 33.37
$ws.Range("Z19").Value = 34.69
$ws.Range("U21").Value = -7.2
$ws.Range("V21").Value = -6.91
$ws.Range("W21").Value = -5.19
$ws.Range("X21").Value = -2.22
$ws.Range("Y21").Value = -1.6
$ws.Range("Z21").Value = -0.97
$ws.Range("U24").Value = 93.52
$ws.Range("V24").Value = 42.67
$ws.Range("W24").Value = 76.26000000000001
$ws.Range("X24").Value = 40.34
$ws.Range("Y24").Value = 33.37
$ws.Range("Z24").Value = 34.69
$ws.Range("U25").Value = 0
$ws.Range("V25").Value = -48.26
$ws.Range("U26").Value = -7.2
$ws.Range("V26").Value = -6.91
$ws.Range("W26").Value = -5.19
$ws.Range("X26").Value = -2.22
$ws.Range("Y26").Value = -1.6
$ws.Range("Z26").Value = -0.97
$ws.Range("U29").Value = 92.83
$ws.Range("V29").Value = 90.26000000000001
$ws.Range("W29").Value = 75.76000000000001
$ws.Range("X29").Value = 40.11
$ws.Range("Y29").Value = 33.18
$ws.Range("Z29").Value = 34.66
$ws.Range("U31").Value = -7.89
$ws.Range("V31").Value = -7.58
$ws.Range("W31").Value = -5.68
$ws.Range("X31").Value = -2.45
$ws.Range("Y31").Value = -1.79
$ws.Range("Z31").Value = -1.01
$ws.Range("U34").Value = 96.47
$ws.Range("V34").Value = 94.98999999999999
$ws.Range("W34").Value = 80.63
$ws.Range("X34").Value = 43.61
$ws.Range("Y34").Value = 35.91
$ws.Range("Z34").Value = 37.54
$ws.Range("U36").Value = -4.24
$ws.Range("V36").Value = -2.85
$ws.Range("W36").Value = -0.8100000000000001
$ws.Range("X36").Value = 1.05
$ws.Range("Y36").Value = 0.93
$ws.Range("Z36").Value = 1.88
$ws.Range("V39").Value = 43.18
$ws.Range("W39").Value = 76.61
$ws.Range("X39").Value = 40.53
$ws.Range("Y39").Value = 33.47
$ws.Range("Z39").Value = 34.66
$ws.Range("U40").Value = 0
$ws.Range("V40").Value = -48.26
$ws.Range("U41").Value = -6.59
$ws.Range("V41").Value = -6.4
$ws.Range("W41").Value = -4.83
$ws.Range("X41").Value = -2.03
$ws.Range("Y41").Value = -1.51
$ws.Range("Z41").Value = -1.01
$ws.Range("U44").Value = 96.23999999999999
$ws.Range("V44").Value = 93.27
$ws.Range("W44").Value = 78.2
$ws.Range("X44").Value = 41.4
$ws.Range("Y44").Value = 34.04
$ws.Range("Z44").Value = 34.88
$ws.Range("U46").Value = -4.48
$ws.Range("V46").Value = -4.57
$ws.Range("W46").Value = -3.25
$ws.Range("X46").Value = -1.16
$ws.Range("Y46").Value = -0.9399999999999999
$ws.Range("Z46").Value = -0.78
$ws.Range("U49").Value = 100.92
$ws.Range("V49").Value = 102.56
$ws.Range("W49").Value = 85.01000000000001
$ws.Range("X49").Value = 44.01
$ws.Range("Y49").Value = 36.06
$ws.Range("Z49").Value = 35.92
$ws.Range("U51").Value = 0.2
$ws.Range("V51").Value = 4.72
$ws.Range("W51").Value = 3.57
$ws.Range("X51").Value = 1.45
$ws.Range("Z51").Value = 0.25
$ws.Range("U54").Value = 99.62
$ws.Range("V54").Value = 97.55
$ws.Range("W54").Value = 81.12
$ws.Range("X54").Value = 43.25
$ws.Range("Y54").Value = 35.54
$ws.Range("Z54").Value = 36.21
$ws.Range("U56").Value = -1.1
$ws.Range("V56").Value = -0.29
$ws.Range("W56").Value = -0.32
$ws.Range("X56").Value = 0.6899999999999999
$ws.Range("Y56").Value = 0.57
$ws.Range("Z56").Value = 0.54
$ws.Range("U59").Value = 104.26
$ws.Range("V59").Value = 100.87
$ws.Range("W59").Value = 84.83
$ws.Range("X59").Value = 44.38
$ws.Range("Y59").Value = 36.47
$ws.Range("Z59").Value = 37.19
$ws.Range("U61").Value = 3.54
$ws.Range("V61").Value = 3.03
$ws.Range("W61").Value = 3.39
$ws.Range("X61").Value = 1.82
$ws.Range("Y61").Value = 1.5
$ws.Range("Z61").Value = 1.52
$ws.Range("U64").Value = 105.24
$ws.Range("V64").Value = 101.71
$ws.Range("W64").Value = 85.45999999999999
$ws.Range("X64").Value = 44.57
$ws.Range("Y64").Value = 36.7
$ws.Range("Z64").Value = 37.42
$ws.Range("V66").Value = 3.86
$ws.Range("W66").Value = 4.02
$ws.Range("X66").Value = 2.01
$ws.Range("Y66").Value = 1.72
$ws.Range("Z66").Value = 1.76
$ws.Range("U69").Value = 106.02
$ws.Range("V69").Value = 102.34
$ws.Range("W69").Value = 86.55
$ws.Range("X69").Value = 45.28
$ws.Range("Y69").Value = 37.21
$ws.Range("Z69").Value = 37.94
$ws.Range("U71").Value = 5.3
$ws.Range("V71").Value = 4.5
$ws.Range("W71").Value = 5.11
$ws.Range("X71").Value = 2.72
$ws.Range("Z71").Value = 2.28
$ws.Range("U74").Value = 103.73
$ws.Range("V74").Value = 100.76
$ws.Range("W74").Value = 84.56999999999999
$ws.Range("X74").Value = 44.33
$ws.Range("Y74").Value = 36.39
$ws.Range("Z74").Value = 37.07
$ws.Range("V76").Value = 2.92
$ws.Range("W76").Value = 3.13
$ws.Range("X76").Value = 1.77
$ws.Range("Y76").Value = 1.42
$ws.Range("Z76").Value = 1.41
$ws.Range("U79").Value = 104.15
$ws.Range("V79").Value = 101.16
$ws.Range("W79").Value = 84.88
$ws.Range("X79").Value = 44.41
$ws.Range("Y79").Value = 36.54
$ws.Range("Z79").Value = 37.19
$ws.Range("U81").Value = 3.43
$ws.Range("V81").Value = 3.32
$ws.Range("W81").Value = 3.44
$ws.Range("X81").Value = 1.85
$ws.Range("Y81").Value = 1.56
$ws.Range("Z81").Value = 1.52
$ws.Range("U84").Value = 95.56
$ws.Range("V84").Value = 93.72
$ws.Range("W84").Value = 78.08
$ws.Range("X84").Value = 42.73
$ws.Range("Y84").Value = 35.12
$ws.Range("Z84").Value = 35.67
$ws.Range("U86").Value = -5.16
$ws.Range("V86").Value = -4.12
$ws.Range("W86").Value = -3.36
$ws.Range("X86").Value = 0.17
$ws.Range("Y86").Value = 0.14
$ws.Range("U89").Value = 92.83
$ws.Range("V89").Value = 42
$ws.Range("W89").Value = 75.76000000000001
$ws.Range("X89").Value = 40.11
$ws.Range("Y89").Value = 33.18
$ws.Range("Z89").Value = 34.66
$ws.Range("U90").Value = 0
$ws.Range("V90").Value = -48.26
$ws.Range("U91").Value = -7.89
$ws.Range("V91").Value = -7.58
$ws.Range("W91").Value = -5.68
$ws.Range("X91").Value = -2.45
$ws.Range("Y91").Value = -1.79
$ws.Range("Z91").Value = -1.01
